$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 214 (existing data at row 214 and below shifts down by 3).
$ws.Range("A214:T216").EntireRow.Insert()

# New row 214: Larry Ann, Especial
$ws.Range("A214").Value = 3
$ws.Range("B214").Value = "Femacal de La Calera"
$ws.Range("C214").Value = "Coquimbo"
$ws.Range("D214").Value = 44931
$ws.Range("E214").Value = 5
$ws.Range("F214").Value = "Fruta"
$ws.Range("G214").Value = 100103
$ws.Range("H214").Value = "Frutos de hueso (carozo)"
$ws.Range("I214").Value = 100103002
$ws.Range("J214").Value = "Ciruela"
$ws.Range("K214").Value = "Larry Ann"
$ws.Range("L214").Value = "Especial"
$ws.Range("M214").Value = 70
$ws.Range("N214").Value = 18000
$ws.Range("O214").Value = 18000
$ws.Range("P214").Value = 18000
$ws.Range("Q214").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R214").Value = "Región de O'Higgins"
$ws.Range("S214").Value = 1200
$ws.Range("T214").Value = 15

# New row 215: Larry Ann, Primera
$ws.Range("A215").Value = 3
$ws.Range("B215").Value = "Femacal de La Calera"
$ws.Range("C215").Value = "Coquimbo"
$ws.Range("D215").Value = 44931
$ws.Range("E215").Value = 5
$ws.Range("F215").Value = "Fruta"
$ws.Range("G215").Value = 100103
$ws.Range("H215").Value = "Frutos de hueso (carozo)"
$ws.Range("I215").Value = 100103002
$ws.Range("J215").Value = "Ciruela"
$ws.Range("K215").Value = "Larry Ann"
$ws.Range("L215").Value = "Primera"
$ws.Range("M215").Value = 68
$ws.Range("N215").Value = 15000
$ws.Range("O215").Value = 15000
$ws.Range("P215").Value = 15000
$ws.Range("Q215").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R215").Value = "Región de O'Higgins"
$ws.Range("S215").Value = 1000
$ws.Range("T215").Value = 15

# New row 216: Larry Ann, Segunda
$ws.Range("A216").Value = 3
$ws.Range("B216").Value = "Femacal de La Calera"
$ws.Range("C216").Value = "Coquimbo"
$ws.Range("D216").Value = 44931
$ws.Range("E216").Value = 5
$ws.Range("F216").Value = "Fruta"
$ws.Range("G216").Value = 100103
$ws.Range("H216").Value = "Frutos de hueso (carozo)"
$ws.Range("I216").Value = 100103002
$ws.Range("J216").Value = "Ciruela"
$ws.Range("K216").Value = "Larry Ann"
$ws.Range("L216").Value = "Segunda"
$ws.Range("M216").Value = 52
$ws.Range("N216").Value = 12000
$ws.Range("O216").Value = 12000
$ws.Range("P216").Value = 12000
$ws.Range("Q216").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R216").Value = "Región de O'Higgins"
$ws.Range("S216").Value = 800
$ws.Range("T216").Value = 15
